# Regenerate orders with updated distance/size labels.
# The experiment's distance/size condition codes changed:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# These tokens appear as standalone lookup values (e.g. "D64", "S30") and
# embedded inside composite condition/filename strings (e.g.
# "Face12_D80_S30", "Face12_D80_S30_l.png"), so a blanket text
# find-and-replace over the whole used range updates every occurrence
# consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

$rng.Replace("D64", "D69")
$rng.Replace("D80", "D86")
$rng.Replace("D51", "D55")
$rng.Replace("S30", "S31")
